# Updates cryptocurrency price/volume figures (columns D and E) for rows 2-51
# to reflect the latest data refresh, per GitHub Actions automation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Range, [string]$Text)
    # Force the cell to plain text so numeric-looking strings (e.g. "12.10",
    # "0.00001101") are preserved exactly instead of being parsed as numbers,
    # then reset the style back to Normal so no residual number format sticks.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextCellValue $ws.Range("D2") "29.093.13"
Set-TextCellValue $ws.Range("E2") "  -3.81%  "
Set-TextCellValue $ws.Range("D3") "1.965.61"
Set-TextCellValue $ws.Range("E3") "  -5.73%  "
Set-TextCellValue $ws.Range("E4") "  -0.04%  "
Set-TextCellValue $ws.Range("D5") "328.28"
Set-TextCellValue $ws.Range("E5") "  -3.45%  "
Set-TextCellValue $ws.Range("D6") "1.005"
Set-TextCellValue $ws.Range("E6") "  +0.01%  "
Set-TextCellValue $ws.Range("D7") "0.5003"
Set-TextCellValue $ws.Range("E7") "  -5.11%  "
Set-TextCellValue $ws.Range("D8") "0.4221"
Set-TextCellValue $ws.Range("E8") "  -3.61%  "
Set-TextCellValue $ws.Range("D9") "52.88"
Set-TextCellValue $ws.Range("E9") "  -3.66%  "
Set-TextCellValue $ws.Range("D10") "0.09180"
Set-TextCellValue $ws.Range("E10") "  -1.66%  "
Set-TextCellValue $ws.Range("D11") "1.101"
Set-TextCellValue $ws.Range("E11") "  -6.22%  "
Set-TextCellValue $ws.Range("D12") "23.02"
Set-TextCellValue $ws.Range("E12") "  -6.00%  "
Set-TextCellValue $ws.Range("D13") "2.012.16"
Set-TextCellValue $ws.Range("E13") "  -3.79%  "
Set-TextCellValue $ws.Range("D14") "7.869"
Set-TextCellValue $ws.Range("E14") "  -7.16%  "
Set-TextCellValue $ws.Range("D15") "6.442"
Set-TextCellValue $ws.Range("E15") "  -6.00%  "
Set-TextCellValue $ws.Range("E16") "  +0.13%  "
Set-TextCellValue $ws.Range("D17") "0.00001101"
Set-TextCellValue $ws.Range("E17") "  -4.67%  "
Set-TextCellValue $ws.Range("D18") "91.53"
Set-TextCellValue $ws.Range("E18") "  -9.85%  "
Set-TextCellValue $ws.Range("D19") "0.06712"
Set-TextCellValue $ws.Range("E19") "  +0.23%  "
Set-TextCellValue $ws.Range("D20") "19.28"
Set-TextCellValue $ws.Range("E20") "  -8.23%  "
Set-TextCellValue $ws.Range("E21") "  +0.14%  "
Set-TextCellValue $ws.Range("D22") "5.972"
Set-TextCellValue $ws.Range("E22") "  -5.07%  "
Set-TextCellValue $ws.Range("D23") "29.127.75"
Set-TextCellValue $ws.Range("E23") "  -3.76%  "
Set-TextCellValue $ws.Range("D24") "12.10"
Set-TextCellValue $ws.Range("E24") "  -2.29%  "
Set-TextCellValue $ws.Range("D25") "2.286"
Set-TextCellValue $ws.Range("E25") "  -1.58%  "
Set-TextCellValue $ws.Range("D26") "2.245.72"
Set-TextCellValue $ws.Range("E26") "  -3.50%  "
Set-TextCellValue $ws.Range("D27") "156.50"
Set-TextCellValue $ws.Range("E27") "  -3.71%  "
Set-TextCellValue $ws.Range("D28") "20.61"
Set-TextCellValue $ws.Range("E28") "  -5.26%  "
Set-TextCellValue $ws.Range("D29") "6.239"
Set-TextCellValue $ws.Range("E29") "  -8.53%  "
Set-TextCellValue $ws.Range("E30") "  -8.64%  "
Set-TextCellValue $ws.Range("D31") "126.59"
Set-TextCellValue $ws.Range("E31") "  -5.22%  "
Set-TextCellValue $ws.Range("D32") "1.047"
Set-TextCellValue $ws.Range("E32") "  -6.91%  "
Set-TextCellValue $ws.Range("D33") "0.09861"
Set-TextCellValue $ws.Range("E33") "  -5.84%  "
Set-TextCellValue $ws.Range("E34") "  -7.86%  "
Set-TextCellValue $ws.Range("E35") "  -7.31%  "
Set-TextCellValue $ws.Range("D36") "3.679"
Set-TextCellValue $ws.Range("E36") "  -6.01%  "
Set-TextCellValue $ws.Range("D37") "0.02436"
Set-TextCellValue $ws.Range("E37") "  -6.61%  "
Set-TextCellValue $ws.Range("D38") "9.047"
Set-TextCellValue $ws.Range("E38") "  -8.45%  "
Set-TextCellValue $ws.Range("E39") "  -2.95%  "
Set-TextCellValue $ws.Range("D40") "0.06363"
Set-TextCellValue $ws.Range("E40") "  -5.62%  "
Set-TextCellValue $ws.Range("D41") "0.6453"
Set-TextCellValue $ws.Range("E41") "  -6.94%  "
Set-TextCellValue $ws.Range("D42") "11.43"
Set-TextCellValue $ws.Range("E42") "  -8.85%  "
Set-TextCellValue $ws.Range("D43") "0.1991"
Set-TextCellValue $ws.Range("E43") "  -9.50%  "
Set-TextCellValue $ws.Range("E44") "  +0.11%  "
Set-TextCellValue $ws.Range("D45") "0.6253"
Set-TextCellValue $ws.Range("E45") "  -7.12%  "
Set-TextCellValue $ws.Range("D46") "13.38"
Set-TextCellValue $ws.Range("E46") "  -6.63%  "
Set-TextCellValue $ws.Range("E47") "  -7.56%  "
Set-TextCellValue $ws.Range("D48") "1.288"
Set-TextCellValue $ws.Range("E48") "  +0.09%  "
Set-TextCellValue $ws.Range("E49") "  -4.42%  "
Set-TextCellValue $ws.Range("D50") "0.00000000332"
Set-TextCellValue $ws.Range("E50") "  -3.23%  "
Set-TextCellValue $ws.Range("D51") "0.06987"
